$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose "Sending cluster" was "ECs"; the six remaining
# rows (previously 5-10, for sending clusters FAPs/MuSCs) shift up to 2-7.
$ws.Rows("2:4").Delete()

# The underlying TPM values were recomputed (ECs no longer contributes to the
# normalization), so refresh every data cell in the six remaining rows with the
# updated figures.
$data = @(
    ,("FAPs", "Sfrp1", "Fzd6", "ECs", 3, 1, 17.23456, 51.70368000000001, 0.8703823568377641, 0.870382356837764, 3, 1, 12.28101533333333, 36.843046, 0.959552102275422, 0.959552102275422, 211.6578956232534, 1904.92106060928, 0.835177220287113, 0.8351772202871129)
    ,("FAPs", "Sfrp1", "Fzd6", "FAPs", 3, 1, 17.23456, 51.70368000000001, 0.8703823568377641, 0.870382356837764, 2, 0.6666666666666666, 0.09168666666666665, 0.27506, 0.007163750827004844, 0.007163750827004845, 1.580179357866667, 14.2216142208, 0.006235202328606958, 0.006235202328606958)
    ,("FAPs", "Sfrp1", "Fzd6", "MuSCs", 3, 1, 17.23456, 51.70368000000001, 0.8703823568377641, 0.870382356837764, 3, 1, 0.4259936666666666, 1.277981, 0.03328414689757318, 0.03328414689757318, 7.341813407786667, 66.07632067008, 0.02896993422204409, 0.02896993422204409)
    ,("MuSCs", "Sfrp1", "Fzd6", "ECs", 3, 1, 2.566576666666667, 7.69973, 0.129617643162236, 0.129617643162236, 3, 1, 12.28101533333333, 36.843046, 0.959552102275422, 0.959552102275422, 31.52016739750889, 283.68150657758, 0.124374881988309, 0.124374881988309)
    ,("MuSCs", "Sfrp1", "Fzd6", "FAPs", 3, 1, 2.566576666666667, 7.69973, 0.129617643162236, 0.129617643162236, 2, 0.6666666666666666, 0.09168666666666665, 0.27506, 0.007163750827004844, 0.007163750827004845, 0.2353208593111111, 2.1178877338, 0.0009285484983978867, 0.0009285484983978866)
    ,("MuSCs", "Sfrp1", "Fzd6", "MuSCs", 3, 1, 2.566576666666667, 7.69973, 0.129617643162236, 0.129617643162236, 3, 1, 0.4259936666666666, 1.277981, 0.03328414689757318, 0.03328414689757318, 1.093345405014444, 9.840108645129998, 0.004314212675529084, 0.004314212675529083)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
